# Error Calculations and Plots
#
# The source data table simulates "missing" values by leaving a cell
# blank (present in the sheet as an empty inline string) rather than by
# omitting the cell entirely. Two data rows ("RM 232" and "SC 92") are
# removed outright, which shifts every following row up by one/two, and a
# handful of cells in columns C and E swap which values are known vs.
# missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two rows that were dropped from the source table. Deleting
# the lower-numbered row first (28 = "SC 92") keeps row 26's identity
# ("RM 232") stable for the second delete.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# After the deletions the remaining rows have shifted up so that the
# original "SC 5" row (27) is now row 26, "SC 101" (29) is now row 27,
# and so on through the end of the table ("SC 232" now row 33).

# --- Column E (D) missing-value swaps among the RM rows ---
$ws.Range("E6").Value = -5.7
$ws.Range("E8").ClearContents()
$ws.Range("E18").Value = -8.5
$ws.Range("E20").ClearContents()
$ws.Range("E23").Value = -7
$ws.Range("E25").ClearContents()

# --- Column C (B) / E (D) missing-value swaps among the (shifted) SC rows ---
$ws.Range("C27").Value = 10
$ws.Range("C28").ClearContents()
$ws.Range("C29").ClearContents()
$ws.Range("C30").Value = 11.4
$ws.Range("E30").Value = -5.7
$ws.Range("C32").ClearContents()
